$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.913903365559928
$ws.Range("D2").Value = 4.292552934790345
$ws.Range("E2").Value = 11.2888838682746
$ws.Range("F2").Value = 61.7386625164971
$ws.Range("G2").Value = 3.773113911216053
$ws.Range("J2").Value = 10.87173636623256
$ws.Range("K2").Value = 24.14971029846336

$ws.Range("C3").Value = 4.762165049485147
$ws.Range("D3").Value = 4.29370355723789
$ws.Range("E3").Value = 11.31808361841776
$ws.Range("F3").Value = 60.98766605690877
$ws.Range("G3").Value = 3.778714272135998
$ws.Range("J3").Value = 10.87562178676379
$ws.Range("K3").Value = 23.92483556898414

$ws.Range("C4").Value = 4.668374986915694
$ws.Range("D4").Value = 4.29538028734332
$ws.Range("E4").Value = 11.33897390350713
$ws.Range("F4").Value = 60.53458786347476
$ws.Range("G4").Value = 3.782321379759144
$ws.Range("J4").Value = 10.87999669866327
$ws.Range("K4").Value = 23.79534979399848

$ws.Range("C5").Value = 4.630064456868946
$ws.Range("D5").Value = 4.296304498848992
$ws.Range("E5").Value = 11.34823008507262
$ws.Range("F5").Value = 60.35211590577624
$ws.Range("G5").Value = 3.783833879416628
$ws.Range("J5").Value = 10.88227836002967
$ws.Range("K5").Value = 23.74479965216622

$ws.Range("C6").Value = 4.623699571301194
$ws.Range("D6").Value = 4.296472424315297
$ws.Range("E6").Value = 11.34981190988912
$ws.Range("F6").Value = 60.3219510830892
$ws.Range("G6").Value = 3.784087606411622
$ws.Range("J6").Value = 10.88268731452908
$ws.Range("K6").Value = 23.73654126107861

$ws.Range("C7").Value = 4.667858593326319
$ws.Range("D7").Value = 4.295391780147788
$ws.Range("E7").Value = 11.33909572836044
$ws.Range("F7").Value = 60.53211805168224
$ws.Range("G7").Value = 3.782341605195726
$ws.Range("J7").Value = 10.88002545202582
$ws.Range("K7").Value = 23.79465901183365

$ws.Range("C8").Value = 4.861751608528723
$ws.Range("D8").Value = 4.292746489348708
$ws.Range("E8").Value = 11.29833642691772
$ws.Range("F8").Value = 61.47813667597208
$ws.Range("G8").Value = 3.775010079539688
$ws.Range("J8").Value = 10.87266241304854
$ws.Range("K8").Value = 24.07042661107546

$ws.Range("C9").Value = 5.23433636741937
$ws.Range("D9").Value = 4.295390067235378
$ws.Range("E9").Value = 11.24196764963731
$ws.Range("F9").Value = 63.39065702344281
$ws.Range("G9").Value = 3.761959936327263
$ws.Range("J9").Value = 10.87406346985668
$ws.Range("K9").Value = 24.67661765160344

$ws.Range("C10").Value = 5.5001095884455
$ws.Range("D10").Value = 4.302282192513761
$ws.Range("E10").Value = 11.21499429981989
$ws.Range("F10").Value = 64.82198555822717
$ws.Range("G10").Value = 3.753167282648948
$ws.Range("J10").Value = 10.88481916034132
$ws.Range("K10").Value = 25.15791422085303

$ws.Range("C11").Value = 5.618680627037248
$ws.Range("D11").Value = 4.306527995359287
$ws.Range("E11").Value = 11.20587284300657
$ws.Range("F11").Value = 65.47691611359097
$ws.Range("G11").Value = 3.749337015099598
$ws.Range("J11").Value = 10.89183587369689
$ws.Range("K11").Value = 25.38376343302624

$ws.Range("C12").Value = 5.663201055956439
$ws.Range("D12").Value = 4.308298182266298
$ws.Range("E12").Value = 11.20287241857629
$ws.Range("F12").Value = 65.72530848794939
$ws.Range("G12").Value = 3.747910742486342
$ws.Range("J12").Value = 10.89479893856412
$ws.Range("K12").Value = 25.47020051685418

$ws.Range("C13").Value = 5.653630342533369
$ws.Range("D13").Value = 4.307909676640003
$ws.Range("E13").Value = 11.20349842434014
$ws.Range("F13").Value = 65.67179824228026
$ws.Range("G13").Value = 3.748216844432742
$ws.Range("J13").Value = 10.89414717338676
$ws.Range("K13").Value = 25.45154542983011

$ws.Range("C14").Value = 5.622351202180395
$ws.Range("D14").Value = 4.30667036005389
$ws.Range("E14").Value = 11.20561689960924
$ws.Range("F14").Value = 65.49734473310852
$ws.Range("G14").Value = 3.749219191695837
$ws.Range("J14").Value = 10.89207351171601
$ws.Range("K14").Value = 25.39085683090199

$ws.Range("C15").Value = 5.603141108743955
$ws.Range("D15").Value = 4.305932469376423
$ws.Range("E15").Value = 11.20697363149771
$ws.Range("F15").Value = 65.39053214269693
$ws.Range("G15").Value = 3.749836299134023
$ws.Range("J15").Value = 10.89084319672784
$ws.Range("K15").Value = 25.35379980083767

$ws.Range("C16").Value = 5.492309867985194
$ws.Range("D16").Value = 4.302027299653273
$ws.Range("E16").Value = 11.21565384988329
$ws.Range("F16").Value = 64.77924698587006
$ws.Range("G16").Value = 3.753420993143052
$ws.Range("J16").Value = 10.88440339296644
$ws.Range("K16").Value = 25.14328657720001

$ws.Range("C17").Value = 5.423688301840195
$ws.Range("D17").Value = 4.299917998353488
$ws.Range("E17").Value = 11.22178605397358
$ws.Range("F17").Value = 64.40510281902571
$ws.Range("G17").Value = 3.755663361028641
$ws.Range("J17").Value = 10.88099722739997
$ws.Range("K17").Value = 25.01585555455432

$ws.Range("C18").Value = 5.38400200786463
$ws.Range("D18").Value = 4.298809159336723
$ws.Range("E18").Value = 11.22560947934067
$ws.Range("F18").Value = 64.19027355771459
$ws.Range("G18").Value = 3.756969083947733
$ws.Range("J18").Value = 10.87923797917323
$ws.Range("K18").Value = 24.94321542871259

$ws.Range("C19").Value = 5.370529081623968
$ws.Range("D19").Value = 4.298451567882662
$ws.Range("E19").Value = 11.22695489342095
$ws.Range("F19").Value = 64.11760437190249
$ws.Range("G19").Value = 3.757413929073175
$ws.Range("J19").Value = 10.87867663851596
$ws.Range("K19").Value = 24.91873561428039

$ws.Range("C20").Value = 5.431015980273344
$ws.Range("D20").Value = 4.300131714704629
$ws.Range("E20").Value = 11.22110259308099
$ws.Range("F20").Value = 64.44489414936278
$ws.Range("G20").Value = 3.755423005631716
$ws.Range("J20").Value = 10.88133912716065
$ws.Range("K20").Value = 25.02935360805838

$ws.Range("C21").Value = 5.631549286696378
$ws.Range("D21").Value = 4.307029949142477
$ws.Range("E21").Value = 11.20498233340163
$ws.Range("F21").Value = 65.54857682358841
$ws.Range("G21").Value = 3.748924123852673
$ws.Range("J21").Value = 10.89267428829802
$ws.Range("K21").Value = 25.40865841162719

$ws.Range("C22").Value = 5.760377782992969
$ws.Range("D22").Value = 4.312486065488604
$ws.Range("E22").Value = 11.19709145040991
$ws.Range("F22").Value = 66.27206301149839
$ws.Range("G22").Value = 3.744817489644417
$ws.Range("J22").Value = 10.90186614377323
$ws.Range("K22").Value = 25.66183923992719

$ws.Range("C23").Value = 5.691837384215992
$ws.Range("D23").Value = 4.309486456041552
$ws.Range("E23").Value = 11.20106073466934
$ws.Range("F23").Value = 65.88578033972554
$ws.Range("G23").Value = 3.746996469349245
$ws.Range("J23").Value = 10.89679692743985
$ws.Range("K23").Value = 25.52625516348441

$ws.Range("C24").Value = 5.427703864889549
$ws.Range("D24").Value = 4.300034770352299
$ws.Range("E24").Value = 11.22141065780359
$ws.Range("F24").Value = 64.42690363174017
$ws.Range("G24").Value = 3.755531618723065
$ws.Range("J24").Value = 10.88118393447841
$ws.Range("K24").Value = 25.02324919566005

$ws.Range("C25").Value = 5.134719453682552
$ws.Range("D25").Value = 4.293817982322104
$ws.Range("E25").Value = 11.25468598578296
$ws.Range("F25").Value = 62.8680229196567
$ws.Range("G25").Value = 3.765349695612366
$ws.Range("J25").Value = 10.8719807876905
$ws.Range("K25").Value = 24.50601796186861
